# RPA datasets push 2023-11-19
# Update the IPO/SPAC tracking table on Sheet1:
#   - Insert a new top row for "하나스팩30호" (Hana SPAC No.30)
#   - Insert a new row for "디에스단석(구,단석산업)" (DS Dansuk) after 키움스팩9호
#   - Every other existing row shifts down by two rows accordingly
#   - The two oldest rows (에이직랜드, 에스와이스틸텍) fall off the bottom of the table
#   - The confirmed offer price ("확정공모가") for 에이텀 is updated from "-" to "18000"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Full target data for rows 2..21 (columns A..F), already shifted/merged per the diff.
$data = @(
    @("하나스팩30호",                 "2023.12.08~12.11", "2,000~2,000",   "-",     14000,  "하나증권"),
    @("키움스팩9호",                  "2023.12.07~12.08", "2,000~2,000",   "-",     6000,   "키움증권"),
    @("디에스단석(구,단석산업)",      "2023.12.05~12.11", "79,000~89,000", "-",     96380,  "KB증권,NH투자증권"),
    @("블루엠텍",                     "2023.11.22~11.28", "15,000~19,000", "-",     21000,  "하나증권,키움증권"),
    @("LS머트리얼즈",                 "2023.11.22~11.28", "4,400~5,500",   "-",     64350,  "키움증권,KB증권,이베스트투자증권,하이투자증권,NH투자증권"),
    @("삼성스팩9호",                  "2023.11.20~11.21", "2,000~2,000",   "-",     20000,  "삼성증권"),
    @("교보스팩15호",                 "2023.11.20~11.21", "2,000~2,000",   "-",     7000,   "교보증권"),
    @("케이엔에스",                   "2023.11.16~11.22", "19,000~22,000", "-",     14250,  "신영증권"),
    @("NH스팩30호",                   "2023.11.15~11.16", "2,000~2,000",   "-",     16000,  "NH투자증권"),
    @("와이바이오로직스",             "2023.11.10~11.16", "9,000~11,000",  "-",     13500,  "유안타증권"),
    @("에이텀",                       "2023.11.09~11.15", "23,000~30,000", "18000", 14950,  "하나증권"),
    @("에이에스텍",                   "2023.11.07~11.13", "21,000~25,000", "28000", 29547,  "미래에셋증권"),
    @("그린리소스",                   "2023.11.03~11.09", "11,000~14,000", "17000", 18040,  "NH투자증권"),
    @("한선엔지니어링",               "2023.11.02~11.08", "5,200~6,000",   "7000",  22100,  "대신증권"),
    @("에코아이",                     "2023.11.01~11.07", "28,500~34,700", "34700", 59251,  "KB증권"),
    @("동인기연(유가)",               "2023.11.01~11.07", "33,000~37,000", "30000", 60654,  "NH투자증권"),
    @("스톰테크",                     "2023.10.31~11.06", "8,000~9,500",   "11000", 26800,  "하이투자증권"),
    @("에코프로머티리얼즈(유가)",     "2023.10.30~11.03", "36,200~44,000", "36200", 524031, "미래에셋증권,NH투자증권,하이투자증권"),
    @("캡스톤파트너스",               "2023.10.26~11.01", "3,200~3,600",   "4000",  6384,   "NH투자증권"),
    @("한국스팩13호",                 "2023.10.25~10.26", "2,000~2,000",   "2000",  8000,   "한국투자증권")
)

$r = 2
foreach ($row in $data) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $r++
}
